# This script reproduces the FlashScore weekly-games update:
#   - A new match (Queretaro vs Santos Laguna, MEXICO - LIGA MX) is inserted as row 3.
#   - The existing row 3 (Racing Montevideo vs Nacional, URUGUAY - PRIMERA DIVISION) is
#     pushed down to row 4, and several of its odds columns are refreshed with new values.
#   - The sheet dimension grows from A1:BD3 to A1:BD4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3; this shifts the old row 3 (and anything below it) down
# by one row, so the old row 3 becomes row 4. Row/column formatting is carried along.
$ws.Rows(3).Insert()

# Column B holds dates written as plain text (e.g. "08/11/2024"). Without forcing a Text
# number format first, Excel would silently reinterpret that string as a date serial value.
$ws.Range("B3:B4").NumberFormat = "@"

# --- Row 3 (new match): MEXICO - LIGA MX, Queretaro vs Santos Laguna ---
$row3 = New-Object 'object[,]' 1,56
$row3[0,0] = 'OGYJXPXN'  # Id
$row3[0,1] = '08/11/2024'  # Date
$row3[0,2] = '22:00'  # Time
$row3[0,3] = 'MEXICO - LIGA MX'  # League
$row3[0,4] = 'Queretaro'  # Home
$row3[0,5] = 'Santos Laguna'  # Away
$row3[0,6] = 2  # Odd_H_FT
$row3[0,7] = 3.6  # Odd_D_FT
$row3[0,8] = 3.6  # Odd_A_FT
$row3[0,9] = 2.63  # Odd_H_HT
$row3[0,10] = 2.2  # Odd_D_HT
$row3[0,11] = 4  # Odd_A_HT
$row3[0,12] = 1.05  # Odd_Over05_FT
$row3[0,13] = 11  # Odd_Under05_FT
$row3[0,14] = 1.29  # Odd_Over15_FT
$row3[0,15] = 3.5  # Odd_Under15_FT
$row3[0,16] = 1.93  # Odd_Over25_FT
$row3[0,17] = 1.93  # Odd_Under25_FT
$row3[0,18] = 1.4  # Odd_Over05_HT
$row3[0,19] = 2.75  # Odd_Under05_HT
$row3[0,20] = 1.75  # Odd_BTTS_Yes
$row3[0,21] = 2  # Odd_BTTS_No
$row3[0,22] = 8  # Odd_CS_1-0
$row3[0,23] = 10  # Odd_CS_2-0
$row3[0,24] = 9  # Odd_CS_2-1
$row3[0,25] = 17  # Odd_CS_3-0
$row3[0,26] = 15  # Odd_CS_3-1
$row3[0,27] = 26  # Odd_CS_3-2
$row3[0,28] = 11  # Odd_CS_0-0
$row3[0,29] = 6.5  # Odd_CS_1-1
$row3[0,30] = 15  # Odd_CS_2-2
$row3[0,31] = 51  # Odd_CS_3-3
$row3[0,32] = 201  # Odd_CS_4-4
$row3[0,33] = 11  # Odd_CS_0-1
$row3[0,34] = 19  # Odd_CS_0-2
$row3[0,35] = 12  # Odd_CS_1-2
$row3[0,36] = 41  # Odd_CS_0-3
$row3[0,37] = 29  # Odd_CS_1-3
$row3[0,38] = 34  # Odd_CS_2-3
$row3[0,39] = 4  # Odd_CS_1-0_HT
$row3[0,40] = 11  # Odd_CS_2-0_HT
$row3[0,41] = 21  # Odd_CS_2-1_HT
$row3[0,42] = 41  # Odd_CS_3-0_HT
$row3[0,43] = 51  # Odd_CS_3-1_HT
$row3[0,44] = 151  # Odd_CS_3-2_HT
$row3[0,45] = 2.75  # Odd_CS_0-0_HT
$row3[0,46] = 8  # Odd_CS_1-1_HT
$row3[0,47] = 51  # Odd_CS_2-2_HT
$row3[0,48] = 5.5  # Odd_CS_0-1_HT
$row3[0,49] = 19  # Odd_CS_0-2_HT
$row3[0,50] = 26  # Odd_CS_1-2_HT
$row3[0,51] = 67  # Odd_CS_0-3_HT
$row3[0,52] = 81  # Odd_CS_1-3_HT
$row3[0,53] = 151  # Odd_CS_2-3_HT
$row3[0,54] = 501  # Odd_CS_3-3_HT
$row3[0,55] = 126  # Odd_CS_4-4_HT
$ws.Range("A3:BD3").Value = $row3

# --- Row 4 (shifted match, odds refreshed): URUGUAY - PRIMERA DIVISION, Racing Montevideo vs Nacional ---
$row4 = New-Object 'object[,]' 1,56
$row4[0,0] = 'QmuqFgzh'  # Id
$row4[0,1] = '08/11/2024'  # Date
$row4[0,2] = '20:30'  # Time
$row4[0,3] = 'URUGUAY - PRIMERA DIVISION'  # League
$row4[0,4] = 'Racing Montevideo'  # Home
$row4[0,5] = 'Nacional'  # Away
$row4[0,6] = 8  # Odd_H_FT
$row4[0,7] = 3.9  # Odd_D_FT
$row4[0,8] = 1.5  # Odd_A_FT
$row4[0,9] = 8  # Odd_H_HT
$row4[0,10] = 2.05  # Odd_D_HT
$row4[0,11] = 2.1  # Odd_A_HT
$row4[0,12] = 1.1  # Odd_Over05_FT
$row4[0,13] = 7  # Odd_Under05_FT
$row4[0,14] = 1.44  # Odd_Over15_FT
$row4[0,15] = 2.63  # Odd_Under15_FT
$row4[0,16] = 2.4  # Odd_Over25_FT
$row4[0,17] = 1.53  # Odd_Under25_FT
$row4[0,18] = 1.53  # Odd_Over05_HT
$row4[0,19] = 2.38  # Odd_Under05_HT
$row4[0,20] = 2.63  # Odd_BTTS_Yes
$row4[0,21] = 1.44  # Odd_BTTS_No
$row4[0,22] = 13  # Odd_CS_1-0
$row4[0,23] = 34  # Odd_CS_2-0
$row4[0,24] = 23  # Odd_CS_2-1
$row4[0,25] = 101  # Odd_CS_3-0
$row4[0,26] = 67  # Odd_CS_3-1
$row4[0,27] = 81  # Odd_CS_3-2
$row4[0,28] = 7  # Odd_CS_0-0
$row4[0,29] = 8  # Odd_CS_1-1
$row4[0,30] = 29  # Odd_CS_2-2
$row4[0,31] = 101  # Odd_CS_3-3
$row4[0,32] = 351  # Odd_CS_4-4
$row4[0,33] = 4.75  # Odd_CS_0-1
$row4[0,34] = 5.5  # Odd_CS_0-2
$row4[0,35] = 9.5  # Odd_CS_1-2
$row4[0,36] = 9.5  # Odd_CS_0-3
$row4[0,37] = 17  # Odd_CS_1-3
$row4[0,38] = 41  # Odd_CS_2-3
$row4[0,39] = 8.5  # Odd_CS_1-0_HT
$row4[0,40] = 41  # Odd_CS_2-0_HT
$row4[0,41] = 51  # Odd_CS_2-1_HT
$row4[0,42] = 201  # Odd_CS_3-0_HT
$row4[0,43] = 301  # Odd_CS_3-1_HT
$row4[0,44] = 301  # Odd_CS_3-2_HT
$row4[0,45] = 2.38  # Odd_CS_0-0_HT
$row4[0,46] = 11  # Odd_CS_1-1_HT
$row4[0,47] = 101  # Odd_CS_2-2_HT
$row4[0,48] = 3.2  # Odd_CS_0-1_HT
$row4[0,49] = 8  # Odd_CS_0-2_HT
$row4[0,50] = 26  # Odd_CS_1-2_HT
$row4[0,51] = 26  # Odd_CS_0-3_HT
$row4[0,52] = 51  # Odd_CS_1-3_HT
$row4[0,53] = 251  # Odd_CS_2-3_HT
$row4[0,54] = 51  # Odd_CS_3-3_HT
$row4[0,55] = 51  # Odd_CS_4-4_HT
$ws.Range("A4:BD4").Value = $row4

Write-Host "Done updating rows 3 and 4"
